# UndoRedoExecuteUndoStateListDiagram.pptx correction:
#  1) Bump the cached "datetimeFigureOut" footer field on every slide
#     layout + the slide master from 20/3/19 -> 21/3/19.
#  2) Merge the two runs that make up "currentStatePointer = 2" / "= 1"
#     textboxes on the slide into a single run (dropping the stray
#     err="1" spellcheck flag left over from the split run).

$p = $ppt.ActivePresentation

# --- 1) Date placeholder fields -------------------------------------------
$newDate = "21/3/19"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.Designs.Item(1).SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Merge "currentStatePointer" runs on slide 1 ------------------------
$slide = $p.Slides.Item(1)

function Merge-CurrentStatePointerRun($shape, $fullText) {
    $tr = $shape.TextFrame.TextRange
    # "currentStatePointer" is exactly the first 19 characters of the run;
    # clearing just that span removes the stray run (with its err="1" flag)
    # while leaving the second run (and its rPr/formatting) untouched.
    $first = $tr.Characters(1, 19)
    $first.Text = ""
    $rest = $shape.TextFrame.TextRange
    $rest.Text = $fullText
}

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "currentStatePointer = 2") {
            Merge-CurrentStatePointerRun $shp "currentStatePointer = 2"
        } elseif ($txt -eq "currentStatePointer = 1") {
            Merge-CurrentStatePointerRun $shp "currentStatePointer = 1"
        }
    }
}
